# Concentration check feature: normalize the "ML" unit suffix to lowercase "ml"
# in the Concentração (column D) values, per the shared-strings diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "1000ml"
$ws.Range("D3").Value = "250ml"
$ws.Range("D4").Value = "500ml"
$ws.Range("D5").Value = "1000ml"
